# Update "想去人数" (want-to-go count) figures that changed between the
# previous scrape and the refreshed one (gh-pages output regenerated at
# commit 456a3b4).
#
# Sheet "展览" (exhibitions):
#   F3  11569 -> 11572
#   F7  11528 -> 11532
#   F10    83 -> 84
#   F12  5721 -> 5722
#
# Sheet "全部类型" (all types, a combined view of every sheet):
#   F5  11569 -> 11572
#   F9  11528 -> 11532
#   F12    83 -> 84
#   F15  5721 -> 5722

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 11572
$wsExhibit.Range("F7").Value = 11532
$wsExhibit.Range("F10").Value = 84
$wsExhibit.Range("F12").Value = 5722

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 11572
$wsAll.Range("F9").Value = 11532
$wsAll.Range("F12").Value = 84
$wsAll.Range("F15").Value = 5722
